$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.520.05'
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.840.65'
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -4.39%  '
$c.ClearFormats()

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9993'
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.52%  '
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '312.92'
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -3.94%  '
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.44%  '
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4225'
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -7.94%  '
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3632'
$c.ClearFormats()

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = 'Dogecoin'
$c.ClearFormats()

$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07207'
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -6.98%  '
$c.ClearFormats()

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = 'Polygon'
$c.ClearFormats()

$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8996'
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -8.21%  '
$c.ClearFormats()

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = 'Solana'
$c.ClearFormats()

$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.57'
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -9.01%  '
$c.ClearFormats()

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = 'WrappedEther'
$c.ClearFormats()

$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.843.43'
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -5.57%  '
$c.ClearFormats()

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = 'Chainlink'
$c.ClearFormats()

$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.580'
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -5.63%  '
$c.ClearFormats()

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = 'Polkadot'
$c.ClearFormats()

$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.303'
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -7.05%  '
$c.ClearFormats()

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = 'TRON'
$c.ClearFormats()

$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.06794'
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -3.38%  '
$c.ClearFormats()

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = 'BinanceUSD'
$c.ClearFormats()

$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.48%  '
$c.ClearFormats()

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = 'Litecoin'
$c.ClearFormats()

$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '77.21'
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -8.97%  '
$c.ClearFormats()

$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = 'ShibaInu'
$c.ClearFormats()

$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008930'
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -6.04%  '
$c.ClearFormats()

$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c.ClearFormats()

$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.9993'
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.42%  '
$c.ClearFormats()

$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = 'Avalanche'
$c.ClearFormats()

$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c.ClearFormats()

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '15.34'
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -8.02%  '
$c.ClearFormats()

$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = 'WrappedBTC'
$c.ClearFormats()

$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '27.507.22'
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -5.38%  '
$c.ClearFormats()

$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = 'Uniswap'
$c.ClearFormats()

$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.936'
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -7.90%  '
$c.ClearFormats()

$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = 'Cosmos'
$c.ClearFormats()

$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.64'
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -3.70%  '
$c.ClearFormats()

$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = 'WrappedliquidstakedEther2.0'
$c.ClearFormats()

$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.086.62'
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -4.06%  '
$c.ClearFormats()

$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = 'Toncoin'
$c.ClearFormats()

$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.042'
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -0.61%  '
$c.ClearFormats()

$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = 'Monero'
$c.ClearFormats()

$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '151.02'
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -4.51%  '
$c.ClearFormats()

$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = 'EthereumClassic'
$c.ClearFormats()

$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c.ClearFormats()

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.20'
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -4.30%  '
$c.ClearFormats()

$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c.ClearFormats()

$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.273'
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -6.11%  '
$c.ClearFormats()

$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = 'BitcoinCash'
$c.ClearFormats()

$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '110.60'
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -6.03%  '
$c.ClearFormats()

$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c.ClearFormats()

$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.688'
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -7.82%  '
$c.ClearFormats()

$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = 'Stellar'
$c.ClearFormats()

$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08862'
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -4.97%  '
$c.ClearFormats()

$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c.ClearFormats()

$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7757'
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -9.99%  '
$c.ClearFormats()

$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c.ClearFormats()

$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.503'
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -11.65%  '
$c.ClearFormats()

$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'HuobiToken'
$c.ClearFormats()

$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.858'
$c.ClearFormats()

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -5.22%  '
$c.ClearFormats()

$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c.ClearFormats()

$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.080'
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -13.25%  '
$c.ClearFormats()

$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = 'Frax'
$c.ClearFormats()

$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.9986'
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c.ClearFormats()

$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = 'Hedera'
$c.ClearFormats()

$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.05385'
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -5.38%  '
$c.ClearFormats()

$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = 'TrustWalletToken'
$c.ClearFormats()

$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.092'
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -5.04%  '
$c.ClearFormats()

$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = 'VeChain'
$c.ClearFormats()

$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01932'
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -5.54%  '
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.952'
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -4.57%  '
$c.ClearFormats()

$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = 'TheSandbox'
$c.ClearFormats()

$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.5031'
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -8.67%  '
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.773'
$c.ClearFormats()

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -8.99%  '
$c.ClearFormats()

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.ClearFormats()

$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1635'
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -6.93%  '
$c.ClearFormats()

$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c.ClearFormats()

$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.06609'
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -4.72%  '
$c.ClearFormats()

$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = 'Aptos'
$c.ClearFormats()

$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.216'
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -12.05%  '
$c.ClearFormats()

$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = 'Decentraland'
$c.ClearFormats()

$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4720'
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -8.94%  '
$c.ClearFormats()

$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.ClearFormats()

$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '105.30'
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -4.96%  '
$c.ClearFormats()

$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c.ClearFormats()

$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '10.22'
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -9.17%  '
$c.ClearFormats()

$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = 'PaxDollar'
$c.ClearFormats()

$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.9983'
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c.ClearFormats()

$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = 'NEARProtocol'
$c.ClearFormats()

$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.639'
$c.ClearFormats()

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -6.97%  '
$c.ClearFormats()

$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.ClearFormats()

$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.842'
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -15.81%  '
$c.ClearFormats()

